# Update the "Förändrad" date column (C) for rows 2-15 from
# 2023-10-09 (serial 45208) to 2023-10-13 (serial 45212).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45208) {
        $cell.Value = 45212
    }
}
